$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.295.91"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "1.874.35"
$ws.Range("E3").Value = "  +0.67%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7122"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.62"
$ws.Range("E6").Value = "  +0.99%  "

$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3109"
$ws.Range("E8").Value = "  +1.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07731"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.11"
$ws.Range("E10").Value = "  +1.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08458"
$ws.Range("E11").Value = "  +2.53%  "

$ws.Range("D12").Value = "1.867.12"
$ws.Range("E12").Value = "  +0.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.211"
$ws.Range("E13").Value = "  +0.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7107"
$ws.Range("E14").Value = "  -0.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.33"
$ws.Range("E15").Value = "  +1.29%  "

$ws.Range("D16").Value = "29.301.04"
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008300"
$ws.Range("E17").Value = "  +6.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.979"
$ws.Range("E18").Value = "  +2.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.59"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").Value = "2.127.77"
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("E21").Value = "  +0.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.798"
$ws.Range("E23").Value = "  -1.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1626"
$ws.Range("E25").Value = "  +2.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.13"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.008"
$ws.Range("E27").Value = "  +1.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.51"
$ws.Range("E28").Value = "  +1.95%  "

$ws.Range("E29").Value = "  +0.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.418"
$ws.Range("E30").Value = "  +1.75%  "

$ws.Range("E31").Value = "  +6.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.277"
$ws.Range("E32").Value = "  -4.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05257"
$ws.Range("E33").Value = "  +1.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.920"
$ws.Range("E34").Value = "  +0.46%  "

$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7452"
$ws.Range("E36").Value = "  +2.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.682"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("E38").Value = "  +0.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.725"
$ws.Range("E39").Value = "  +1.13%  "

$ws.Range("D40").Value = "1.162.70"
$ws.Range("E40").Value = "  +0.60%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.359"
$ws.Range("E41").Value = "  +4.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8891"
$ws.Range("E42").Value = "  -1.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.89"
$ws.Range("E43").Value = "  +1.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.78"
$ws.Range("E44").Value = "  +5.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9998"

$ws.Range("D46").Value = "2.023.46"
$ws.Range("E46").Value = "  +0.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.804"
$ws.Range("E47").Value = "  +2.51%  "

$ws.Range("E48").Value = "  -1.66%  "

$ws.Range("E49").Value = "  +3.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.379"
$ws.Range("E50").Value = "  +1.15%  "

$ws.Range("E51").Value = "  +1.62%  "
